{"js": "// Replace each old three-digit-division answer with its new value.\n// The mapping below mirrors the document's table cells in reading order\n// (each old string is unique in the document, so a single-hit search is safe).\nconst replacements = [\n  [\"288\u00f79=32, 0\", \"484\u00f76=80, 4\"],\n  [\"898\u00f78=112, 2\", \"120\u00f79=13, 3\"],\n  [\"424\u00f74=106, 0\", \"930\u00f76=155, 0\"],\n  [\"130\u00f77=18, 4\", \"712\u00f73=237, 1\"],\n  [\"230\u00f77=32, 6\", \"421\u00f79=46, 7\"],\n  [\"355\u00f72=177, 1\", \"113\u00f74=28, 1\"],\n  [\"361\u00f78=45, 1\", \"466\u00f79=51, 7\"],\n  [\"269\u00f74=67, 1\", \"855\u00f77=122, 1\"],\n  [\"552\u00f74=138, 0\", \"847\u00f79=94, 1\"],\n  [\"584\u00f78=73, 0\", \"726\u00f75=145, 1\"],\n  [\"310\u00f78=38, 6\", \"855\u00f79=95, 0\"],\n  [\"530\u00f73=176, 2\", \"300\u00f72=150, 0\"],\n  [\"674\u00f76=112, 2\", \"154\u00f74=38, 2\"],\n  [\"245\u00f73=81, 2\", \"750\u00f76=125, 0\"],\n  [\"699\u00f74=174, 3\", \"397\u00f75=79, 2\"],\n  [\"640\u00f73=213, 1\", \"829\u00f75=165, 4\"],\n  [\"926\u00f72=463, 0\", \"474\u00f78=59, 2\"],\n  [\"363\u00f79=40, 3\", \"462\u00f79=51, 3\"],\n  [\"286\u00f74=71, 2\", \"235\u00f76=39, 1\"],\n  [\"285\u00f76=47, 3\", \"800\u00f75=160, 0\"],\n  [\"670\u00f79=74, 4\", \"868\u00f73=289, 1\"],\n  [\"182\u00f72=91, 0\", \"487\u00f79=54, 1\"],\n  [\"454\u00f72=227, 0\", \"187\u00f77=26, 5\"],\n  [\"131\u00f77=18, 5\", \"429\u00f78=53, 5\"],\n  [\"912\u00f75=182, 2\", \"471\u00f79=52, 3\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each old three-digit-division answer with its new value.\n# Each \"Old\" string is unique within the document, so Find.Execute with\n# wdReplaceAll (2) safely touches exactly the one matching table cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"288\u00f79=32, 0\"; New = \"484\u00f76=80, 4\" }\n    @{ Old = \"898\u00f78=112, 2\"; New = \"120\u00f79=13, 3\" }\n    @{ Old = \"424\u00f74=106, 0\"; New = \"930\u00f76=155, 0\" }\n    @{ Old = \"130\u00f77=18, 4\"; New = \"712\u00f73=237, 1\" }\n    @{ Old = \"230\u00f77=32, 6\"; New = \"421\u00f79=46, 7\" }\n    @{ Old = \"355\u00f72=177, 1\"; New = \"113\u00f74=28, 1\" }\n    @{ Old = \"361\u00f78=45, 1\"; New = \"466\u00f79=51, 7\" }\n    @{ Old = \"269\u00f74=67, 1\"; New = \"855\u00f77=122, 1\" }\n    @{ Old = \"552\u00f74=138, 0\"; New = \"847\u00f79=94, 1\" }\n    @{ Old = \"584\u00f78=73, 0\"; New = \"726\u00f75=145, 1\" }\n    @{ Old = \"310\u00f78=38, 6\"; New = \"855\u00f79=95, 0\" }\n    @{ Old = \"530\u00f73=176, 2\"; New = \"300\u00f72=150, 0\" }\n    @{ Old = \"674\u00f76=112, 2\"; New = \"154\u00f74=38, 2\" }\n    @{ Old = \"245\u00f73=81, 2\"; New = \"750\u00f76=125, 0\" }\n    @{ Old = \"699\u00f74=174, 3\"; New = \"397\u00f75=79, 2\" }\n    @{ Old = \"640\u00f73=213, 1\"; New = \"829\u00f75=165, 4\" }\n    @{ Old = \"926\u00f72=463, 0\"; New = \"474\u00f78=59, 2\" }\n    @{ Old = \"363\u00f79=40, 3\"; New = \"462\u00f79=51, 3\" }\n    @{ Old = \"286\u00f74=71, 2\"; New = \"235\u00f76=39, 1\" }\n    @{ Old = \"285\u00f76=47, 3\"; New = \"800\u00f75=160, 0\" }\n    @{ Old = \"670\u00f79=74, 4\"; New = \"868\u00f73=289, 1\" }\n    @{ Old = \"182\u00f72=91, 0\"; New = \"487\u00f79=54, 1\" }\n    @{ Old = \"454\u00f72=227, 0\"; New = \"187\u00f77=26, 5\" }\n    @{ Old = \"131\u00f77=18, 5\"; New = \"429\u00f78=53, 5\" }\n    @{ Old = \"912\u00f75=182, 2\"; New = \"471\u00f79=52, 3\" }\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap(=wdFindContinue), Format,\n    #         ReplaceWith, Replace(=wdReplaceAll))\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
